# Auto-generated COM-interop script applying the ZBP_06_home_office.xlsx update
# (new weekly columns 13.-19. 12. 2021 and 20.-26. 12. 2021, "aktualizace" date bump)
$wb = $excel.ActiveWorkbook

# ===================== Sheet "data" =====================
$ws1 = $wb.Worksheets.Item("data")

# Header row 1: new week-label columns BT1/BU1, formatted like BS1
$ws1.Range("BS1").Copy() | Out-Null
$ws1.Range("BT1").PasteSpecial(-4122) | Out-Null
$ws1.Range("BT1").Value = "13.–19. 12. 2021"
$ws1.Range("BU1").PasteSpecial(-4122) | Out-Null
$ws1.Range("BU1").Value = "20.–26. 12. 2021"
$excel.CutCopyMode = $false

# Data rows 2-77: BT = week 13.-19.12.2021, BU = week 20.-26.12.2021
$ws1.Cells.Item(2, 72).Value = 0.75
$ws1.Cells.Item(2, 73).Value = 0.57
$ws1.Cells.Item(3, 72).Value = 0.07000000000000001
$ws1.Cells.Item(3, 73).Value = 0.13
$ws1.Cells.Item(4, 72).Value = 0.08
$ws1.Cells.Item(4, 73).Value = 0.08
$ws1.Cells.Item(5, 72).Value = 0.1
$ws1.Cells.Item(5, 73).Value = 0.22
$ws1.Cells.Item(6, 72).Value = 0.74
$ws1.Cells.Item(6, 73).Value = 0.61
$ws1.Cells.Item(7, 72).Value = 0.06
$ws1.Cells.Item(7, 73).Value = 0.09
$ws1.Cells.Item(8, 72).Value = 0.09
$ws1.Cells.Item(8, 73).Value = 0.09
$ws1.Cells.Item(9, 72).Value = 0.11
$ws1.Cells.Item(9, 73).Value = 0.21
$ws1.Cells.Item(10, 72).Value = 0.47
$ws1.Cells.Item(10, 73).Value = 0.25
$ws1.Cells.Item(11, 72).Value = 0.15
$ws1.Cells.Item(11, 73).Value = 0.14
$ws1.Cells.Item(12, 72).Value = 0.32
$ws1.Cells.Item(12, 73).Value = 0.35
$ws1.Cells.Item(13, 72).Value = 0.06
$ws1.Cells.Item(13, 73).Value = 0.26
$ws1.Cells.Item(14, 72).Value = 0.83
$ws1.Cells.Item(14, 73).Value = 0.6
$ws1.Cells.Item(15, 72).Value = 0.04
$ws1.Cells.Item(15, 73).Value = 0.14
$ws1.Cells.Item(16, 72).Value = 0.02
$ws1.Cells.Item(16, 73).Value = 0.03
$ws1.Cells.Item(17, 72).Value = 0.11
$ws1.Cells.Item(17, 73).Value = 0.23
$ws1.Cells.Item(18, 72).Value = 0.84
$ws1.Cells.Item(18, 73).Value = 0.6899999999999999
$ws1.Cells.Item(19, 72).Value = 0.05
$ws1.Cells.Item(19, 73).Value = 0.11
$ws1.Cells.Item(20, 72).Value = 0.03
$ws1.Cells.Item(20, 73).Value = 0.04
$ws1.Cells.Item(21, 72).Value = 0.08
$ws1.Cells.Item(21, 73).Value = 0.16
$ws1.Cells.Item(22, 72).Value = 0.75
$ws1.Cells.Item(22, 73).Value = 0.43
$ws1.Cells.Item(23, 72).Value = 0.13
$ws1.Cells.Item(23, 73).Value = 0.34
$ws1.Cells.Item(24, 72).Value = 0.01
$ws1.Cells.Item(24, 73).Value = 0.03
$ws1.Cells.Item(25, 72).Value = 0.11
$ws1.Cells.Item(25, 73).Value = 0.2
$ws1.Cells.Item(26, 72).Value = 0.76
$ws1.Cells.Item(26, 73).Value = 0.58
$ws1.Cells.Item(27, 72).Value = 0.05
$ws1.Cells.Item(27, 73).Value = 0.12
$ws1.Cells.Item(28, 72).Value = 0.06
$ws1.Cells.Item(28, 73).Value = 0.07000000000000001
$ws1.Cells.Item(29, 72).Value = 0.13
$ws1.Cells.Item(29, 73).Value = 0.23
$ws1.Cells.Item(30, 72).Value = 0.82
$ws1.Cells.Item(30, 73).Value = 0.63
$ws1.Cells.Item(31, 72).Value = 0.07000000000000001
$ws1.Cells.Item(31, 73).Value = 0.14
$ws1.Cells.Item(32, 72).Value = 0.05
$ws1.Cells.Item(32, 73).Value = 0.05
$ws1.Cells.Item(33, 72).Value = 0.06
$ws1.Cells.Item(33, 73).Value = 0.18
$ws1.Cells.Item(34, 72).Value = 0.63
$ws1.Cells.Item(34, 73).Value = 0.47
$ws1.Cells.Item(35, 72).Value = 0.09
$ws1.Cells.Item(35, 73).Value = 0.08
$ws1.Cells.Item(36, 72).Value = 0.15
$ws1.Cells.Item(36, 73).Value = 0.14
$ws1.Cells.Item(37, 72).Value = 0.13
$ws1.Cells.Item(37, 73).Value = 0.31
$ws1.Cells.Item(38, 72).Value = 0.68
$ws1.Cells.Item(38, 73).Value = 0.51
$ws1.Cells.Item(39, 72).Value = 0.09
$ws1.Cells.Item(39, 73).Value = 0.16
$ws1.Cells.Item(40, 72).Value = 0.14
$ws1.Cells.Item(40, 73).Value = 0.14
$ws1.Cells.Item(41, 72).Value = 0.09
$ws1.Cells.Item(41, 73).Value = 0.19
$ws1.Cells.Item(42, 72).Value = 0.71
$ws1.Cells.Item(42, 73).Value = 0.5600000000000001
$ws1.Cells.Item(43, 72).Value = 0.12
$ws1.Cells.Item(43, 73).Value = 0.13
$ws1.Cells.Item(44, 72).Value = 0.11
$ws1.Cells.Item(44, 73).Value = 0.1
$ws1.Cells.Item(45, 72).Value = 0.06
$ws1.Cells.Item(45, 73).Value = 0.21
$ws1.Cells.Item(46, 72).Value = 0.8100000000000001
$ws1.Cells.Item(46, 73).Value = 0.61
$ws1.Cells.Item(47, 72).Value = 0.04
$ws1.Cells.Item(47, 73).Value = 0.13
$ws1.Cells.Item(48, 72).Value = 0.06
$ws1.Cells.Item(48, 73).Value = 0.07000000000000001
$ws1.Cells.Item(49, 72).Value = 0.09
$ws1.Cells.Item(49, 73).Value = 0.19
$ws1.Cells.Item(50, 72).Value = 0.64
$ws1.Cells.Item(50, 73).Value = 0.46
$ws1.Cells.Item(51, 72).Value = 0.07000000000000001
$ws1.Cells.Item(51, 73).Value = 0.12
$ws1.Cells.Item(52, 72).Value = 0.09
$ws1.Cells.Item(52, 73).Value = 0.09
$ws1.Cells.Item(53, 72).Value = 0.2
$ws1.Cells.Item(53, 73).Value = 0.33
$ws1.Cells.Item(54, 72).Value = 0.73
$ws1.Cells.Item(54, 73).Value = 0.58
$ws1.Cells.Item(55, 72).Value = 0.05
$ws1.Cells.Item(55, 73).Value = 0.1
$ws1.Cells.Item(56, 72).Value = 0.08
$ws1.Cells.Item(56, 73).Value = 0.09
$ws1.Cells.Item(57, 72).Value = 0.14
$ws1.Cells.Item(57, 73).Value = 0.23
$ws1.Cells.Item(58, 72).Value = 0.8
$ws1.Cells.Item(58, 73).Value = 0.57
$ws1.Cells.Item(59, 72).Value = 0.08
$ws1.Cells.Item(59, 73).Value = 0.13
$ws1.Cells.Item(60, 72).Value = 0.04
$ws1.Cells.Item(60, 73).Value = 0.08
$ws1.Cells.Item(61, 72).Value = 0.08
$ws1.Cells.Item(61, 73).Value = 0.22
$ws1.Cells.Item(62, 72).Value = 0.75
$ws1.Cells.Item(62, 73).Value = 0.5600000000000001
$ws1.Cells.Item(63, 72).Value = 0.08
$ws1.Cells.Item(63, 73).Value = 0.15
$ws1.Cells.Item(64, 72).Value = 0.09
$ws1.Cells.Item(64, 73).Value = 0.07000000000000001
$ws1.Cells.Item(65, 72).Value = 0.08
$ws1.Cells.Item(65, 73).Value = 0.22
$ws1.Cells.Item(66, 72).Value = 0.78
$ws1.Cells.Item(66, 73).Value = 0.6
$ws1.Cells.Item(67, 72).Value = 0.03
$ws1.Cells.Item(67, 73).Value = 0.08
$ws1.Cells.Item(68, 72).Value = 0.03
$ws1.Cells.Item(68, 73).Value = 0.015
$ws1.Cells.Item(69, 72).Value = 0.16
$ws1.Cells.Item(69, 73).Value = 0.305
$ws1.Cells.Item(70, 72).Value = 0.75
$ws1.Cells.Item(70, 73).Value = 0.6
$ws1.Cells.Item(71, 72).Value = 0.08
$ws1.Cells.Item(71, 73).Value = 0.15
$ws1.Cells.Item(72, 72).Value = 0.1
$ws1.Cells.Item(72, 73).Value = 0.08
$ws1.Cells.Item(73, 72).Value = 0.07000000000000001
$ws1.Cells.Item(73, 73).Value = 0.17
$ws1.Cells.Item(74, 72).Value = 0.72
$ws1.Cells.Item(74, 73).Value = 0.47
$ws1.Cells.Item(75, 72).Value = 0.11
$ws1.Cells.Item(75, 73).Value = 0.18
$ws1.Cells.Item(76, 72).Value = 0.13
$ws1.Cells.Item(76, 73).Value = 0.19
$ws1.Cells.Item(77, 72).Value = 0.04
$ws1.Cells.Item(77, 73).Value = 0.16

# Row 78 label: bump "aktualizace" date
$ws1.Range("A78").Value = "Život během pandemie, Home office, % respondentů celkově a ve skupinách, aktualizace 6. 1. 2022"

# ===================== Sheet "pocetR" =====================
$ws2 = $wb.Worksheets.Item("pocetR")

# Header row 1: new week-label columns BS1/BT1, formatted like BR1
$ws2.Range("BR1").Copy() | Out-Null
$ws2.Range("BS1").PasteSpecial(-4122) | Out-Null
$ws2.Range("BS1").Value = "13.–19. 12. 2021"
$ws2.Range("BT1").PasteSpecial(-4122) | Out-Null
$ws2.Range("BT1").Value = "20.–26. 12. 2021"
$excel.CutCopyMode = $false

# Data rows 2-20: BS = week 13.-19.12.2021, BT = week 20.-26.12.2021 (sample sizes)
$ws2.Cells.Item(2, 71).Value = 994
$ws2.Cells.Item(2, 72).Value = 994
$ws2.Cells.Item(3, 71).Value = 273
$ws2.Cells.Item(3, 72).Value = 273
$ws2.Cells.Item(4, 71).Value = 87
$ws2.Cells.Item(4, 72).Value = 87
$ws2.Cells.Item(5, 71).Value = 275
$ws2.Cells.Item(5, 72).Value = 275
$ws2.Cells.Item(6, 71).Value = 149
$ws2.Cells.Item(6, 72).Value = 149
$ws2.Cells.Item(7, 71).Value = 79
$ws2.Cells.Item(7, 72).Value = 79
$ws2.Cells.Item(8, 71).Value = 487
$ws2.Cells.Item(8, 72).Value = 487
$ws2.Cells.Item(9, 71).Value = 264
$ws2.Cells.Item(9, 72).Value = 264
$ws2.Cells.Item(10, 71).Value = 112
$ws2.Cells.Item(10, 72).Value = 112
$ws2.Cells.Item(11, 71).Value = 131
$ws2.Cells.Item(11, 72).Value = 131
$ws2.Cells.Item(12, 71).Value = 264
$ws2.Cells.Item(12, 72).Value = 264
$ws2.Cells.Item(13, 71).Value = 552
$ws2.Cells.Item(13, 72).Value = 552
$ws2.Cells.Item(14, 71).Value = 178
$ws2.Cells.Item(14, 72).Value = 178
$ws2.Cells.Item(15, 71).Value = 402
$ws2.Cells.Item(15, 72).Value = 402
$ws2.Cells.Item(16, 71).Value = 218
$ws2.Cells.Item(16, 72).Value = 218
$ws2.Cells.Item(17, 71).Value = 374
$ws2.Cells.Item(17, 72).Value = 374
$ws2.Cells.Item(18, 71).Value = 381
$ws2.Cells.Item(18, 72).Value = 381
$ws2.Cells.Item(19, 71).Value = 380
$ws2.Cells.Item(19, 72).Value = 380
$ws2.Cells.Item(20, 71).Value = 233
$ws2.Cells.Item(20, 72).Value = 233

# Row 21: label bump + blank placeholder cells trailing the row, formatted like BR21
$ws2.Range("BR21").Copy() | Out-Null
$ws2.Range("BS21").PasteSpecial(-4122) | Out-Null
$ws2.Range("BT21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws2.Range("A21").Value = "Život během pandemie, Home office, velikost dotázaného souboru celkově a ve skupinách, aktualizace 6. 1. 2022"

